$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (style/border/font) from column BF into the new BG, BH, BI columns
$ws.Range("BF1:BF80").Copy($ws.Range("BG1:BG80"))
$ws.Range("BF1:BF80").Copy($ws.Range("BH1:BH80"))
$ws.Range("BF1:BF80").Copy($ws.Range("BI1:BI80"))

# Now overwrite the copied values with the actual new-quarter data
$ws.Cells.Item(1, 59).Value = "31/12/2023"
$ws.Cells.Item(1, 60).Value = "31/03/2024"
$ws.Cells.Item(1, 61).Value = "30/06/2024"
$ws.Cells.Item(2, 59).Value = 3076137.984
$ws.Cells.Item(2, 60).Value = 2936739.072
$ws.Cells.Item(2, 61).Value = 2962825.984
$ws.Cells.Item(3, 59).Value = 138882
$ws.Cells.Item(3, 60).Value = 107087
$ws.Cells.Item(3, 61).Value = 161231.008
$ws.Cells.Item(4, 59).Value = 17600
$ws.Cells.Item(4, 60).Value = 14321
$ws.Cells.Item(4, 61).Value = 78182
$ws.Cells.Item(5, 59).Value = 5180
$ws.Cells.Item(5, 60).Value = 4781
$ws.Cells.Item(5, 61).Value = 4868
$ws.Cells.Item(6, 59).Value = 32930
$ws.Cells.Item(6, 60).Value = 26745
$ws.Cells.Item(6, 61).Value = 30600
$ws.Cells.Item(7, 59).Value = 0
$ws.Cells.Item(7, 60).Value = 0
$ws.Cells.Item(7, 61).Value = 0
$ws.Cells.Item(8, 59).Value = 0
$ws.Cells.Item(8, 60).Value = 0
$ws.Cells.Item(8, 61).Value = 0
$ws.Cells.Item(9, 59).Value = 20300
$ws.Cells.Item(9, 60).Value = 17452
$ws.Cells.Item(9, 61).Value = 12961
$ws.Cells.Item(10, 59).Value = 1620
$ws.Cells.Item(10, 60).Value = 1024
$ws.Cells.Item(10, 61).Value = 1028
$ws.Cells.Item(11, 59).Value = 61252
$ws.Cells.Item(11, 60).Value = 42764
$ws.Cells.Item(11, 61).Value = 33592
$ws.Cells.Item(12, 59).Value = 315255.008
$ws.Cells.Item(12, 60).Value = 226976
$ws.Cells.Item(12, 61).Value = 217099.008
$ws.Cells.Item(13, 59).Value = 0
$ws.Cells.Item(13, 60).Value = 0
$ws.Cells.Item(13, 61).Value = 0
$ws.Cells.Item(14, 59).Value = 0
$ws.Cells.Item(14, 60).Value = 0
$ws.Cells.Item(14, 61).Value = 0
$ws.Cells.Item(15, 59).Value = 312
$ws.Cells.Item(15, 60).Value = 315
$ws.Cells.Item(15, 61).Value = 21
$ws.Cells.Item(16, 59).Value = 0
$ws.Cells.Item(16, 60).Value = 0
$ws.Cells.Item(16, 61).Value = 0
$ws.Cells.Item(17, 59).Value = 0
$ws.Cells.Item(17, 60).Value = 0
$ws.Cells.Item(17, 61).Value = 0
$ws.Cells.Item(18, 59).Value = 0
$ws.Cells.Item(18, 60).Value = 0
$ws.Cells.Item(18, 61).Value = 0
$ws.Cells.Item(19, 59).Value = 4799
$ws.Cells.Item(19, 60).Value = 5102
$ws.Cells.Item(19, 61).Value = 3529
$ws.Cells.Item(20, 59).Value = 0
$ws.Cells.Item(20, 60).Value = 0
$ws.Cells.Item(20, 61).Value = 0
$ws.Cells.Item(21, 59).Value = 0
$ws.Cells.Item(21, 60).Value = 0
$ws.Cells.Item(21, 61).Value = 0
$ws.Cells.Item(22, 59).Value = 0
$ws.Cells.Item(22, 60).Value = 0
$ws.Cells.Item(22, 61).Value = 0
$ws.Cells.Item(23, 59).Value = 2622000.896
$ws.Cells.Item(23, 60).Value = 2602675.968
$ws.Cells.Item(23, 61).Value = 2584496.128
$ws.Cells.Item(24, 59).Value = 0
$ws.Cells.Item(24, 60).Value = 0
$ws.Cells.Item(24, 61).Value = 0
$ws.Cells.Item(25, 59).Value = 0
$ws.Cells.Item(25, 60).Value = 0
$ws.Cells.Item(25, 61).Value = 0
$ws.Cells.Item(26, 59).Value = 3076137.984
$ws.Cells.Item(26, 60).Value = 2936739.072
$ws.Cells.Item(26, 61).Value = 2962825.984
$ws.Cells.Item(27, 59).Value = 437913.984
$ws.Cells.Item(27, 60).Value = 473134.016
$ws.Cells.Item(27, 61).Value = 506687.008
$ws.Cells.Item(28, 59).Value = 7220
$ws.Cells.Item(28, 60).Value = 8766
$ws.Cells.Item(28, 61).Value = 6101
$ws.Cells.Item(29, 59).Value = 63750
$ws.Cells.Item(29, 60).Value = 46843
$ws.Cells.Item(29, 61).Value = 45795
$ws.Cells.Item(30, 59).Value = 15033
$ws.Cells.Item(30, 60).Value = 9511
$ws.Cells.Item(30, 61).Value = 6807
$ws.Cells.Item(31, 59).Value = 172754
$ws.Cells.Item(31, 60).Value = 225503.008
$ws.Cells.Item(31, 61).Value = 253178
$ws.Cells.Item(32, 59).Value = 0
$ws.Cells.Item(32, 60).Value = 0
$ws.Cells.Item(32, 61).Value = 0
$ws.Cells.Item(33, 59).Value = 0
$ws.Cells.Item(33, 60).Value = 0
$ws.Cells.Item(33, 61).Value = 0
$ws.Cells.Item(34, 59).Value = 179156.992
$ws.Cells.Item(34, 60).Value = 182511.008
$ws.Cells.Item(34, 61).Value = 194806
$ws.Cells.Item(35, 59).Value = 0
$ws.Cells.Item(35, 60).Value = 0
$ws.Cells.Item(35, 61).Value = 0
$ws.Cells.Item(36, 59).Value = 0
$ws.Cells.Item(36, 60).Value = 0
$ws.Cells.Item(36, 61).Value = 0
$ws.Cells.Item(37, 59).Value = 1705203.968
$ws.Cells.Item(37, 60).Value = 1592500.992
$ws.Cells.Item(37, 61).Value = 1571192.064
$ws.Cells.Item(38, 59).Value = 901872
$ws.Cells.Item(38, 60).Value = 876820.992
$ws.Cells.Item(38, 61).Value = 876849.9840000001
$ws.Cells.Item(39, 59).Value = 0
$ws.Cells.Item(39, 60).Value = 0
$ws.Cells.Item(39, 61).Value = 0
$ws.Cells.Item(40, 59).Value = 417232
$ws.Cells.Item(40, 60).Value = 407500
$ws.Cells.Item(40, 61).Value = 428267.008
$ws.Cells.Item(41, 59).Value = 0
$ws.Cells.Item(41, 60).Value = 0
$ws.Cells.Item(41, 61).Value = 0
$ws.Cells.Item(42, 59).Value = 0
$ws.Cells.Item(42, 60).Value = 0
$ws.Cells.Item(42, 61).Value = 0
$ws.Cells.Item(43, 59).Value = 386100
$ws.Cells.Item(43, 60).Value = 308180
$ws.Cells.Item(43, 61).Value = 266075.008
$ws.Cells.Item(44, 59).Value = 0
$ws.Cells.Item(44, 60).Value = 0
$ws.Cells.Item(44, 61).Value = 0
$ws.Cells.Item(45, 59).Value = 0
$ws.Cells.Item(45, 60).Value = 0
$ws.Cells.Item(45, 61).Value = 0
$ws.Cells.Item(46, 59).Value = 0
$ws.Cells.Item(46, 60).Value = 0
$ws.Cells.Item(46, 61).Value = 0
$ws.Cells.Item(47, 59).Value = 933020.032
$ws.Cells.Item(47, 60).Value = 871104
$ws.Cells.Item(47, 61).Value = 884947.008
$ws.Cells.Item(48, 59).Value = 4128636.928
$ws.Cells.Item(48, 60).Value = 4128636.928
$ws.Cells.Item(48, 61).Value = 4128636.928
$ws.Cells.Item(49, 59).Value = 1
$ws.Cells.Item(49, 60).Value = 1
$ws.Cells.Item(49, 61).Value = 1
$ws.Cells.Item(50, 59).Value = 0
$ws.Cells.Item(50, 60).Value = 0
$ws.Cells.Item(50, 61).Value = 0
$ws.Cells.Item(51, 59).Value = 0
$ws.Cells.Item(51, 60).Value = 0
$ws.Cells.Item(51, 61).Value = 0
$ws.Cells.Item(52, 59).Value = -3195618.048
$ws.Cells.Item(52, 60).Value = -3257533.952
$ws.Cells.Item(52, 61).Value = -3243691.008
$ws.Cells.Item(53, 59).Value = 0
$ws.Cells.Item(53, 60).Value = 0
$ws.Cells.Item(53, 61).Value = 0
$ws.Cells.Item(54, 59).Value = 0
$ws.Cells.Item(54, 60).Value = 0
$ws.Cells.Item(54, 61).Value = 0
$ws.Cells.Item(55, 59).Value = 0
$ws.Cells.Item(55, 60).Value = 0
$ws.Cells.Item(55, 61).Value = 0
$ws.Cells.Item(56, 59).Value = 0
$ws.Cells.Item(56, 60).Value = 0
$ws.Cells.Item(56, 61).Value = 0
$ws.Cells.Item(59, 59).Value = 50584
$ws.Cells.Item(59, 60).Value = 42587
$ws.Cells.Item(59, 61).Value = 56022
$ws.Cells.Item(60, 59).Value = 104706
$ws.Cells.Item(60, 60).Value = -57538
$ws.Cells.Item(60, 61).Value = -59112
$ws.Cells.Item(61, 59).Value = 155289.984
$ws.Cells.Item(61, 60).Value = -14951
$ws.Cells.Item(61, 61).Value = -3090
$ws.Cells.Item(62, 59).Value = 0
$ws.Cells.Item(62, 60).Value = 0
$ws.Cells.Item(62, 61).Value = 0
$ws.Cells.Item(63, 59).Value = -17934
$ws.Cells.Item(63, 60).Value = -14421
$ws.Cells.Item(63, 61).Value = -11615
$ws.Cells.Item(64, 59).Value = 0
$ws.Cells.Item(64, 60).Value = 0
$ws.Cells.Item(64, 61).Value = 0
$ws.Cells.Item(65, 59).Value = 0
$ws.Cells.Item(65, 60).Value = 0
$ws.Cells.Item(65, 61).Value = 0
$ws.Cells.Item(66, 59).Value = 71126
$ws.Cells.Item(66, 60).Value = -1058
$ws.Cells.Item(66, 61).Value = 62834
$ws.Cells.Item(67, 59).Value = 6878
$ws.Cells.Item(67, 60).Value = 0
$ws.Cells.Item(67, 61).Value = 0
$ws.Cells.Item(68, 59).Value = -31035
$ws.Cells.Item(68, 60).Value = -33814
$ws.Cells.Item(68, 61).Value = -35947
$ws.Cells.Item(69, 59).Value = 757
$ws.Cells.Item(69, 60).Value = 512
$ws.Cells.Item(69, 61).Value = 1378
$ws.Cells.Item(70, 59).Value = -31792.008
$ws.Cells.Item(70, 60).Value = -34326
$ws.Cells.Item(70, 61).Value = -37325
$ws.Cells.Item(74, 59).Value = 184324.992
$ws.Cells.Item(74, 60).Value = -64244
$ws.Cells.Item(74, 61).Value = 12182
$ws.Cells.Item(75, 59).Value = -2193
$ws.Cells.Item(75, 60).Value = -1791
$ws.Cells.Item(75, 61).Value = -16551
$ws.Cells.Item(76, 59).Value = -56276
$ws.Cells.Item(76, 60).Value = 4119
$ws.Cells.Item(76, 61).Value = 18212
$ws.Cells.Item(79, 59).Value = 0
$ws.Cells.Item(79, 60).Value = 0
$ws.Cells.Item(79, 61).Value = 0
$ws.Cells.Item(80, 59).Value = 125856
$ws.Cells.Item(80, 60).Value = -61916
$ws.Cells.Item(80, 61).Value = 13843
